$wb = $excel.ActiveWorkbook

# Reference to the existing "Message" sheet, whose content we duplicate.
$msg = $wb.Worksheets.Item("Message")

# ---------------------------------------------------------------------
# 1) New sheet "CypherOutput_Message" -- exact duplicate of "Message"
# ---------------------------------------------------------------------
$cypherMsg = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$cypherMsg.Name = "CypherOutput_Message"
$msg.Range("A1:A10").Copy()
$cypherMsg.Range("A1").PasteSpecial(-4163)

# ---------------------------------------------------------------------
# 2) New sheet "StatOutput" -- small 4-column stats table
# ---------------------------------------------------------------------
$statOut = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$statOut.Name = "StatOutput"
$statOut.Range("A1").Value = "number_of_files"
$statOut.Range("B1").Value = "number_of_sample"
$statOut.Range("C1").Value = "number_of_cases"
$statOut.Range("D1").Value = "number_of_study"

# Values that look numeric ("117", "62", "27", "1") but must be stored as
# TEXT (shared strings), matching the Neo4j driver's string formatting of
# the result set. Writing them as a TEXT() formula and then flattening the
# formula to a static value via Copy + PasteSpecial(xlPasteValues) keeps
# the text type without registering any extra cell style.
$statOut.Range("A2").Formula = "=TEXT(117,""0"")"
$statOut.Range("B2").Formula = "=TEXT(62,""0"")"
$statOut.Range("C2").Formula = "=TEXT(27,""0"")"
$statOut.Range("D2").Formula = "=TEXT(1,""0"")"
$statRng = $statOut.Range("A2:D2")
$statRng.Copy()
$statRng.PasteSpecial(-4163)

# ---------------------------------------------------------------------
# 3) New sheet "StatOutput_Message" -- "Message" sheet content repeated,
#    but with the Cypher query (row 18) replaced by the new stats query.
# ---------------------------------------------------------------------
$statMsg = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$statMsg.Name = "StatOutput_Message"
$msg.Range("A1:A10").Copy()
$statMsg.Range("A1").PasteSpecial(-4163)
$msg.Range("A1:A7").Copy()
$statMsg.Range("A11").PasteSpecial(-4163)
$statMsg.Range("A18").Value = 'MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.sex IN [''Male'']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study'
$msg.Range("A9:A10").Copy()
$statMsg.Range("A19").PasteSpecial(-4163)

# ---------------------------------------------------------------------
# Restore the originally-selected sheet/tab so CypherOutput keeps focus.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("CypherOutput").Activate()

Write-Host "Workbook now has sheets:"
foreach ($s in $wb.Worksheets) {
    Write-Host (" - " + $s.Name)
}
